$d = $word.ActiveDocument

# 1. Merge the "TUẦN ", "2", ": ", "Architecture & Setup Environment" runs
#    into a single run by replacing the whole heading text in one shot.
$d.Content.Find.Execute(
    "TUẦN 2: Architecture & Setup Environment",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "TUẦN 2: Architecture & Setup Environment",
    2
) | Out-Null

# 2. Nudge the table's inner column widths by 1 dxa each (the first
#    column stays at 1391 dxa / 69.55pt).
$t = $d.Tables.Item(1)
$t.Columns.Item(2).Width = 2342 / 20
$t.Columns.Item(3).Width = 2161 / 20
$t.Columns.Item(4).Width = 1613 / 20
$t.Columns.Item(5).Width = 2465 / 20
